$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66 (shifts existing rows 66-130 down to 67-131,
# matching the "Fruta / hortaliza, semanal" weekly data refresh).
$ws.Rows.Item(66).Insert()

$newRow = 66

$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 45240
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112001
$ws.Cells.Item($newRow, 7).Value = "Berenjena"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 60
$ws.Cells.Item($newRow, 11).Value = 12000
$ws.Cells.Item($newRow, 12).Value = 12000
$ws.Cells.Item($newRow, 13).Value = 12000
$ws.Cells.Item($newRow, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 200
$ws.Cells.Item($newRow, 17).Value = 60
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
